$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Experience -----------------------------------------------------
$ws.Range("B2").Value = 35
$ws.Range("D2").Value = "at least 10 years is acceptable"

# --- Row 3: Price -> Technical Proposal ------------------------------------
$ws.Range("A3").Value = "Technical Proposal"

# --- Row 4: Technical Proposal -> Timeline ----------------------------------
$ws.Range("A4").Value = "Timeline"
$ws.Range("B4").Value = 15

# --- Row 5: Timeline -> References ------------------------------------------
$ws.Range("A5").Value = "References"
$ws.Range("B5").Value = 10
$ws.Range("D5").Value = "a high score will be given for the offer to provide references upon contact, no comment is unacceptable"

# --- Row 6: References -> Quality and OHS Policies (drop C6 entirely) ------
$ws.Rows.Item(6).Delete()
$ws.Range("D6").Value = "alignment with a standard is expected"
$ws.Range("A6").Value = "Quality and OHS Policies"
$ws.Range("B6").Value = 10

$ws.Range("A6").WrapText = $true
$ws.Range("A6").VerticalAlignment = -4108
$ws.Range("B6").WrapText = $true
$ws.Range("B6").VerticalAlignment = -4108
$ws.Range("D6").WrapText = $true
$ws.Range("D6").VerticalAlignment = -4108

# --- Row heights -------------------------------------------------------------
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).RowHeight = 100.8
$ws.Rows.Item(6).RowHeight = 43.2

# --- Selection ----------------------------------------------------------------
$ws.Range("A7").Select()
